# Auto-generated Excel COM-interop script
# Applies numeric cell updates across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# matching the authoritative diff (market-price / profit recalculation refresh).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (58 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1027
$ws.Range("J4").Value = 2739.5
$ws.Range("L4").Value = 2739.5
$ws.Range("N4").Value = -2967.5
$ws.Range("H40").Value = 2844.4443
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 2887.5
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 2887.5
$ws.Range("M40").Value = -2325
$ws.Range("N40").Value = -3237.5
$ws.Range("H41").Value = 67141.39999999999
$ws.Range("I41").Value = 449
$ws.Range("J41").Value = 77401.766
$ws.Range("K41").Value = 449
$ws.Range("L41").Value = 77401.766
$ws.Range("M41").Value = -9
$ws.Range("N41").Value = -78281.766
$ws.Range("H64").Value = 4323.8237
$ws.Range("J64").Value = 4625.375
$ws.Range("L64").Value = 4625.375
$ws.Range("N64").Value = -5121.375
$ws.Range("H67").Value = 4323.8237
$ws.Range("J67").Value = 4625.375
$ws.Range("L67").Value = 4625.375
$ws.Range("N67").Value = -6341.375
$ws.Range("H98").Value = 1358.931
$ws.Range("I98").Value = 1358.931
$ws.Range("K98").Value = 1358.931
$ws.Range("M98").Value = 139.069
$ws.Range("H107").Value = 280.07693
$ws.Range("I107").Value = 376
$ws.Range("J107").Value = 64.25
$ws.Range("K107").Value = 376
$ws.Range("L107").Value = 64.25
$ws.Range("M107").Value = 1544
$ws.Range("N107").Value = -3904.25
$ws.Range("H108").Value = 90000
$ws.Range("J108").Value = 90000
$ws.Range("L108").Value = 90000
$ws.Range("N108").Value = -97680
$ws.Range("I113").Value = 250003490
$ws.Range("J113").Value = 5140.143
$ws.Range("K113").Value = 250003490
$ws.Range("L113").Value = 5140.143
$ws.Range("M113").Value = -250000236
$ws.Range("N113").Value = -11648.143
$ws.Range("H118").Value = 700.1
$ws.Range("I118").Value = 378
$ws.Range("J118").Value = 1988.5
$ws.Range("K118").Value = 1134
$ws.Range("L118").Value = 5965.5
$ws.Range("M118").Value = 523
$ws.Range("N118").Value = -9279.5
$ws.Range("H122").Value = 1358.931
$ws.Range("I122").Value = 1358.931
$ws.Range("K122").Value = 4076.793
$ws.Range("M122").Value = -1626.793

# ---- Sheet: ARM (32 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1218.2632
$ws.Range("I2").Value = 1256.7333
$ws.Range("J2").Value = 1074
$ws.Range("K2").Value = 1256.7333
$ws.Range("L2").Value = 1074
$ws.Range("M2").Value = -1143.7333
$ws.Range("N2").Value = -1300
$ws.Range("H32").Value = 2663.6885
$ws.Range("I32").Value = 1578.5686
$ws.Range("J32").Value = 8197.799999999999
$ws.Range("K32").Value = 1578.5686
$ws.Range("L32").Value = 8197.799999999999
$ws.Range("M32").Value = -1291.5686
$ws.Range("N32").Value = -8771.799999999999
$ws.Range("H45").Value = 1872.3636
$ws.Range("I45").Value = 1399.5555
$ws.Range("K45").Value = 1399.5555
$ws.Range("M45").Value = -1022.5555
$ws.Range("H110").Value = 64492.625
$ws.Range("I110").Value = 43786.082
$ws.Range("J110").Value = 126612.25
$ws.Range("K110").Value = 43786.082
$ws.Range("L110").Value = 126612.25
$ws.Range("M110").Value = -41741.082
$ws.Range("N110").Value = -130702.25
$ws.Range("H116").Value = 1218.2632
$ws.Range("I116").Value = 1256.7333
$ws.Range("J116").Value = 1074
$ws.Range("K116").Value = 1256.7333
$ws.Range("L116").Value = 1074
$ws.Range("M116").Value = 1037.2667
$ws.Range("N116").Value = -5662

# ---- Sheet: BSM (21 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1218.2632
$ws.Range("I3").Value = 1256.7333
$ws.Range("J3").Value = 1074
$ws.Range("K3").Value = 1256.7333
$ws.Range("L3").Value = 1074
$ws.Range("M3").Value = -1142.7333
$ws.Range("N3").Value = -1302
$ws.Range("H105").Value = 1727.4706
$ws.Range("I105").Value = 1621.3077
$ws.Range("J105").Value = 2072.5
$ws.Range("K105").Value = 1621.3077
$ws.Range("L105").Value = 2072.5
$ws.Range("M105").Value = 125.6922999999999
$ws.Range("N105").Value = -5566.5
$ws.Range("H134").Value = 2590.2979
$ws.Range("I134").Value = 2475.3
$ws.Range("J134").Value = 3247.4285
$ws.Range("K134").Value = 7425.900000000001
$ws.Range("L134").Value = 9742.2855
$ws.Range("M134").Value = -4890.900000000001
$ws.Range("N134").Value = -14812.2855

# ---- Sheet: CRP (27 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1763.25
$ws.Range("I16").Value = 1901.6
$ws.Range("K16").Value = 1901.6
$ws.Range("M16").Value = -1614.6
$ws.Range("H105").Value = 1766.3334
$ws.Range("I105").Value = 1766.3334
$ws.Range("K105").Value = 1766.3334
$ws.Range("M105").Value = -19.33339999999998
$ws.Range("H107").Value = 2034.3529
$ws.Range("J107").Value = 3630.2856
$ws.Range("L107").Value = 3630.2856
$ws.Range("N107").Value = -7470.2856
$ws.Range("H113").Value = 1763.25
$ws.Range("I113").Value = 1901.6
$ws.Range("K113").Value = 1901.6
$ws.Range("M113").Value = 268.4000000000001
$ws.Range("H122").Value = 2818
$ws.Range("I122").Value = 2800.375
$ws.Range("J122").Value = 2865
$ws.Range("K122").Value = 8401.125
$ws.Range("L122").Value = 8595
$ws.Range("M122").Value = -5951.125
$ws.Range("N122").Value = -13495
$ws.Range("H132").Value = 2773.4707
$ws.Range("I132").Value = 2564.2856
$ws.Range("K132").Value = 7692.8568
$ws.Range("M132").Value = -5162.8568

# ---- Sheet: CUL (33 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1208.7858
$ws.Range("I8").Value = 1208.7858
$ws.Range("K8").Value = 3626.3574
$ws.Range("M8").Value = -3487.3574
$ws.Range("H10").Value = 560.5263
$ws.Range("I10").Value = 194.18182
$ws.Range("J10").Value = 1064.25
$ws.Range("K10").Value = 582.5454599999999
$ws.Range("L10").Value = 3192.75
$ws.Range("M10").Value = -443.5454599999999
$ws.Range("N10").Value = -3470.75
$ws.Range("H33").Value = 83.69231000000001
$ws.Range("I33").Value = 85.90000000000001
$ws.Range("J33").Value = 76.333336
$ws.Range("K33").Value = 515.4000000000001
$ws.Range("L33").Value = 458.000016
$ws.Range("M33").Value = -232.4000000000001
$ws.Range("N33").Value = -1024.000016
$ws.Range("H74").Value = 4500
$ws.Range("J74").Value = 8000
$ws.Range("L74").Value = 24000
$ws.Range("N74").Value = -26122
$ws.Range("H77").Value = 4500
$ws.Range("J77").Value = 8000
$ws.Range("L77").Value = 72000
$ws.Range("N77").Value = -82608
$ws.Range("H113").Value = 592.3611
$ws.Range("I113").Value = 325.07144
$ws.Range("J113").Value = 762.4545000000001
$ws.Range("K113").Value = 975.21432
$ws.Range("L113").Value = 2287.3635
$ws.Range("M113").Value = 1194.78568
$ws.Range("N113").Value = -6627.3635

# ---- Sheet: GSM (30 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8034.8945
$ws.Range("I70").Value = 7643.6
$ws.Range("J70").Value = 9502.25
$ws.Range("K70").Value = 7643.6
$ws.Range("L70").Value = 9502.25
$ws.Range("M70").Value = -7373.6
$ws.Range("N70").Value = -10042.25
$ws.Range("H73").Value = 8034.8945
$ws.Range("I73").Value = 7643.6
$ws.Range("J73").Value = 9502.25
$ws.Range("K73").Value = 7643.6
$ws.Range("L73").Value = 9502.25
$ws.Range("M73").Value = -6707.6
$ws.Range("N73").Value = -11374.25
$ws.Range("H80").Value = 3967.5
$ws.Range("I80").Value = 4506.1113
$ws.Range("K80").Value = 4506.1113
$ws.Range("M80").Value = -3508.1113
$ws.Range("H83").Value = 3967.5
$ws.Range("I83").Value = 4506.1113
$ws.Range("K83").Value = 22530.5565
$ws.Range("M83").Value = -17538.5565
$ws.Range("H107").Value = 2443.3845
$ws.Range("I107").Value = 1678.1666
$ws.Range("K107").Value = 1678.1666
$ws.Range("M107").Value = 241.8334
$ws.Range("H132").Value = 3883.4
$ws.Range("I132").Value = 3509.238
$ws.Range("K132").Value = 10527.714
$ws.Range("M132").Value = -7997.714

# ---- Sheet: LTW (28 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I20").Value = 994.2
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 994.2
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = -768.2
$ws.Range("N20").Value = -2452
$ws.Range("H46").Value = 1540.5454
$ws.Range("I46").Value = 1841.1666
$ws.Range("K46").Value = 1841.1666
$ws.Range("M46").Value = -1653.1666
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H55").Value = 360.3913
$ws.Range("I55").Value = 350.0625
$ws.Range("J55").Value = 384
$ws.Range("K55").Value = 350.0625
$ws.Range("L55").Value = 384
$ws.Range("M55").Value = -177.0625
$ws.Range("N55").Value = -730
$ws.Range("H132").Value = 6060.385
$ws.Range("I132").Value = 3533
$ws.Range("J132").Value = 10834.333
$ws.Range("K132").Value = 10599
$ws.Range("L132").Value = 32502.999
$ws.Range("M132").Value = -8069
$ws.Range("N132").Value = -37562.999

# ---- Sheet: WVR (33 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 38132
$ws.Range("I21").Value = 34996
$ws.Range("J21").Value = 40013.6
$ws.Range("K21").Value = 34996
$ws.Range("L21").Value = 40013.6
$ws.Range("M21").Value = -34761
$ws.Range("N21").Value = -40483.6
$ws.Range("H35").Value = 38132
$ws.Range("I35").Value = 34996
$ws.Range("J35").Value = 40013.6
$ws.Range("K35").Value = 34996
$ws.Range("L35").Value = 40013.6
$ws.Range("M35").Value = -34706
$ws.Range("N35").Value = -40593.6
$ws.Range("H62").Value = 2651634
$ws.Range("J62").Value = 8558.799999999999
$ws.Range("L62").Value = 8558.799999999999
$ws.Range("N62").Value = -9806.799999999999
$ws.Range("H65").Value = 2651634
$ws.Range("J65").Value = 8558.799999999999
$ws.Range("L65").Value = 42794
$ws.Range("N65").Value = -49034
$ws.Range("H113").Value = 3776.2307
$ws.Range("I113").Value = 359.1
$ws.Range("K113").Value = 1077.3
$ws.Range("M113").Value = 1092.7
$ws.Range("H122").Value = 2630.077
$ws.Range("I122").Value = 2619.2
$ws.Range("J122").Value = 2666.3333
$ws.Range("K122").Value = 7857.599999999999
$ws.Range("L122").Value = 7998.999899999999
$ws.Range("M122").Value = -5407.599999999999
$ws.Range("N122").Value = -12898.9999

# Total cell updates applied: 262
